$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (target stored width ~15.7109375 / 16.42578125 chars;
# ColumnWidth is pixel-quantized on write, so we pick the nearest reachable value)
$ws.Columns.Item(1).ColumnWidth = 14.833333333333332
$ws.Columns.Item(2).ColumnWidth = 15.666666666666668

$ws.Range("A1").Value = -0.31139146475743473
$ws.Range("B1").Value = 0.31079304472140734
$ws.Range("A2").Value = -0.26468698808880475
$ws.Range("B2").Value = 0.26232734122753243
$ws.Range("A3").Value = -0.15938149766746434
$ws.Range("B3").Value = 0.15877212126923368
$ws.Range("A4").Value = -0.1467721215256681
$ws.Range("B4").Value = 0.1462368772724858
$ws.Range("A5").Value = -0.14023687833205667
$ws.Range("B5").Value = 0.13917605005602773
$ws.Range("A6").Value = -0.066604282504426404
$ws.Range("B6").Value = 0.066549681039887343
$ws.Range("A7").Value = -0.046549682303968609
$ws.Range("B7").Value = 0.046444406114265036
$ws.Range("A8").Value = -0.026444407384714097
$ws.Range("B8").Value = 0.026367098388276666
$ws.Range("A9").Value = -0.020367099491266139
$ws.Range("B9").Value = 0.020302046205150859
$ws.Range("A10").Value = -0.014302047312249044
$ws.Range("B10").Value = 0.014301061033762608
$ws.Range("A11").Value = -0.0098010621227366812
$ws.Range("B11").Value = 0.0097841187572207389
$ws.Range("A12").Value = 0.036930371136976614
$ws.Range("B12").Value = -0.037072044048649833
$ws.Range("A13").Value = 0.043072042953443024
$ws.Range("B13").Value = -0.043138921519107143
$ws.Range("A14").Value = -0.0054339783275558773
$ws.Range("B14").Value = 0.0054328761463331432
$ws.Range("A15").Value = 0.00056712276334103251
$ws.Range("B15").Value = -0.00056742642247264996
$ws.Range("A16").Value = 0.0065674253323737553
$ws.Range("B16").Value = -0.0065717681277801532
$ws.Range("A17").Value = 0.012571767039563753
$ws.Range("B17").Value = -0.012578419747601544
$ws.Range("A18").Value = -0.036106057690041382
$ws.Range("B18").Value = 0.036096068493023381
$ws.Range("A19").Value = -0.027096069538134593
$ws.Range("B19").Value = 0.027012917132651193
$ws.Range("A20").Value = -0.018012918186750326
$ws.Range("B20").Value = 0.01800424824834046
$ws.Range("A21").Value = -0.0090042493036843751
$ws.Range("B21").Value = 0.008999998943906462
$ws.Range("A22").Value = -0.093945844634410136
$ws.Range("B22").Value = 0.093634198710731198
$ws.Range("A23").Value = -0.084634199793451437
$ws.Range("B23").Value = 0.084126507170465992
$ws.Range("A24").Value = -0.042126508675643493
$ws.Range("B24").Value = 0.041999998487141532
$ws.Range("A25").Value = -0.066571768640198314
$ws.Range("B25").Value = 0.066517492472378592
$ws.Range("A26").Value = -0.094256188817503528
$ws.Range("B26").Value = 0.094080469953262025
$ws.Range("A27").Value = -0.088080471049001741
$ws.Range("B27").Value = 0.087489516657375699
$ws.Range("A28").Value = -0.081489517776594411
$ws.Range("B28").Value = 0.081101519632188257
$ws.Range("A29").Value = -0.069101520839888408
$ws.Range("B29").Value = 0.068930891669829464
$ws.Range("A30").Value = -0.048930892984421437
$ws.Range("B30").Value = 0.048640293854009808
$ws.Range("A31").Value = -0.027017970624141796
$ws.Range("B31").Value = 0.027000539227998388
$ws.Range("A32").Value = -0.0060005405682739266
$ws.Range("B32").Value = 0.0059999988449481734
